# Bump on-board assignment iteration: duplicate "iter_07" into a new
# "iter_08" sheet (the new, still-being-calibrated iteration), tidy up
# the zoom level on a few of the older sheets, and normalize which
# cells are still highlighted as "changed" now that iter_07's values
# have been carried forward / confirmed.

$wb = $excel.ActiveWorkbook

# --- 1. A few of the older iteration tabs get re-zoomed to 100% -----------
$normalZoomSheets = @("iter_02", "iter_04", "iter_05")
foreach ($name in $normalZoomSheets) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Activate()
    $excel.ActiveWindow.Zoom = 100
}

# --- 2. Duplicate iter_07 -> iter_08 ---------------------------------------
$iter07 = $wb.Worksheets.Item("iter_07")
$iter07.Copy([Type]::Missing, $iter07)
$iter08 = $wb.Worksheets.Item($iter07.Index + 1)
$iter08.Name = "iter_08"

# --- 3. On the new iter_08 sheet: bring in the next round of calibrated ---
#        values for initial_boarding_penalty / transfer_boarding_penalty.
$iter08.Range("B3").Value = 4
$iter08.Range("B4").Value = 4
$iter08.Range("E3").Value = 4.5
$iter08.Range("E4").Value = 4.5
$iter08.Range("F3").Value = 4
$iter08.Range("F4").Value = 4

# These cells kept the same value as iter_07, so they're no longer
# "freshly changed" -- drop the red/bold highlight back to the normal
# black style (copy the formatting from an already-normal cell, B5).
$iter08.Range("B5").Copy()
$normalizeCells = @("D2", "C3", "C4", "D3", "D4", "G3", "G4")
foreach ($addr in $normalizeCells) {
    $iter08.Range($addr).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# --- 4. iter_07 itself: no longer the freshly-edited tab, rezoom + move ---
#        the selection.
$iter07.Activate()
$excel.ActiveWindow.Zoom = 150
$iter07.Range("A2").Select()

# --- 5. iter_08 becomes the active / selected tab --------------------------
$iter08.Activate()
$excel.ActiveWindow.Zoom = 150
$iter08.Range("A1").Select()
